$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -22.1113
$ws.Range("B4").Value = 5.382900000000003

$ws.Range("A7").Value = -20.11699999999997

$ws.Range("C10").Value = -13.60059999999999

$ws.Range("B12").Value = 4.589899999999999

$ws.Range("C13").Value = -13.23649999999999

$ws.Range("A16").Value = -22.00780000000002

$ws.Range("B18").Value = 6.034299999999995

$ws.Range("B19").Value = 8.499999999999998

$ws.Range("B20").Value = 9.518399999999987

$ws.Range("A28").Value = -21.8723

$ws.Range("A29").Value = -21.32189999999997

$ws.Range("C30").Value = -12.27759999999999

$ws.Range("B31").Value = 4.911899999999998

$ws.Range("A32").Value = -21.14289999999999

$ws.Range("A40").Value = -21.93629999999999
$ws.Range("B40").Value = 6.062400000000007
$ws.Range("C40").Value = -12.44340000000001

$ws.Range("B42").Value = 8.560399999999996

$ws.Range("C44").Value = -13.2163

$ws.Range("B47").Value = 5.4544

$ws.Range("B48").Value = 5.355400000000003

$ws.Range("A52").Value = -22.24059999999999

$ws.Range("A57").Value = -21.87920000000001

$ws.Range("B63").Value = 4.753299999999999

$ws.Range("B64").Value = 5.239200000000003

$ws.Range("A66").Value = -21.44290000000001

$ws.Range("B76").Value = 5.122299999999998

$ws.Range("B81").Value = 5.101400000000003

$ws.Range("B89").Value = 5.109299999999995
$ws.Range("C89").Value = -13.4775

$ws.Range("C91").Value = -12.6223

$ws.Range("B94").Value = 4.781799999999993

$ws.Range("A100").Value = -22.05230000000002

$wb.Save()
